{"js": "// Replace each three-digit-by-one-digit multiplication expression with its\n// updated counterpart, per the commit's regenerated worksheet numbers.\n// Each original expression occurs exactly once in the document, so a\n// search-and-replace keyed on the full \"NNN\u00d7N=\" text is unambiguous.\nconst replacements = [\n  [\"529\u00d75=\", \"683\u00d72=\"],\n  [\"422\u00d77=\", \"669\u00d79=\"],\n  [\"841\u00d75=\", \"324\u00d74=\"],\n  [\"256\u00d73=\", \"911\u00d73=\"],\n  [\"542\u00d73=\", \"423\u00d74=\"],\n  [\"230\u00d79=\", \"564\u00d78=\"],\n  [\"689\u00d74=\", \"269\u00d73=\"],\n  [\"243\u00d73=\", \"389\u00d73=\"],\n  [\"701\u00d75=\", \"968\u00d74=\"],\n  [\"779\u00d78=\", \"825\u00d76=\"],\n  [\"527\u00d75=\", \"302\u00d76=\"],\n  [\"515\u00d78=\", \"441\u00d72=\"],\n  [\"866\u00d72=\", \"501\u00d77=\"],\n  [\"378\u00d75=\", \"625\u00d72=\"],\n  [\"774\u00d74=\", \"475\u00d75=\"],\n  [\"333\u00d76=\", \"844\u00d79=\"],\n  [\"702\u00d76=\", \"878\u00d77=\"],\n  [\"476\u00d74=\", \"321\u00d78=\"],\n  [\"755\u00d75=\", \"137\u00d75=\"],\n  [\"247\u00d75=\", \"415\u00d77=\"],\n  [\"259\u00d74=\", \"742\u00d74=\"],\n  [\"952\u00d77=\", \"675\u00d74=\"],\n  [\"200\u00d72=\", \"498\u00d79=\"],\n  [\"548\u00d78=\", \"115\u00d75=\"],\n  [\"439\u00d75=\", \"656\u00d75=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const result of results.items) {\n    result.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication expression with its\n# updated counterpart, per the commit's regenerated worksheet numbers.\n# Each original expression occurs exactly once in the document, so a\n# Find/Replace keyed on the full \"NNN\u00d7N=\" text is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"529\u00d75=\", \"683\u00d72=\"),\n    @(\"422\u00d77=\", \"669\u00d79=\"),\n    @(\"841\u00d75=\", \"324\u00d74=\"),\n    @(\"256\u00d73=\", \"911\u00d73=\"),\n    @(\"542\u00d73=\", \"423\u00d74=\"),\n    @(\"230\u00d79=\", \"564\u00d78=\"),\n    @(\"689\u00d74=\", \"269\u00d73=\"),\n    @(\"243\u00d73=\", \"389\u00d73=\"),\n    @(\"701\u00d75=\", \"968\u00d74=\"),\n    @(\"779\u00d78=\", \"825\u00d76=\"),\n    @(\"527\u00d75=\", \"302\u00d76=\"),\n    @(\"515\u00d78=\", \"441\u00d72=\"),\n    @(\"866\u00d72=\", \"501\u00d77=\"),\n    @(\"378\u00d75=\", \"625\u00d72=\"),\n    @(\"774\u00d74=\", \"475\u00d75=\"),\n    @(\"333\u00d76=\", \"844\u00d79=\"),\n    @(\"702\u00d76=\", \"878\u00d77=\"),\n    @(\"476\u00d74=\", \"321\u00d78=\"),\n    @(\"755\u00d75=\", \"137\u00d75=\"),\n    @(\"247\u00d75=\", \"415\u00d77=\"),\n    @(\"259\u00d74=\", \"742\u00d74=\"),\n    @(\"952\u00d77=\", \"675\u00d74=\"),\n    @(\"200\u00d72=\", \"498\u00d79=\"),\n    @(\"548\u00d78=\", \"115\u00d75=\"),\n    @(\"439\u00d75=\", \"656\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
